# The "Critical X" rolling-average field labels were renamed to "Peak X"
# (commit: "fix name to be peak power").
#
# Sheet "gc_fields_display" (1st sheet) column B holds the human readable
# display name for each internal field name held in column A. Rows 82-85
# correspond to the four __CalcBestRollingWeightedMean* fields whose
# display names need to change from "Critical ..." to "Peak ...".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # gc_fields_display
$ws2 = $wb.Worksheets.Item(2)   # gc_fields_uom

$ws1.Range("B82").Value = "Peak Heart Rate"
$ws1.Range("B83").Value = "Peak Power"
$ws1.Range("B84").Value = "Peak Speed"
$ws1.Range("B85").Value = "Peak Pace"

# Restore the interactive selection/active-sheet state seen in the edited
# workbook (gc_fields_display becomes the active/selected tab, with
# gc_fields_uom's selection parked on C32).
[void]$ws2.Activate()
[void]$ws2.Range("C32").Select()

[void]$ws1.Activate()
[void]$ws1.Range("B86").Select()
